$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.577.35"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "2.541.08"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  -1.08%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "2.938.51"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.88%  "
$ws.Range("D16").Value = "2.569.55"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.836"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "42.574.69"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").Value = "0.0₃0947"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0794"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.110"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0295"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("D46").Value = "1.958.03"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").Value = "2.792.72"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.191"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.28%  "
